$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(3)

# 1) Resize / reposition the shape
$shp.Left = 637032
$shp.Top = 3265765
$shp.Width = 9637776
$shp.Height = 3693319

$tr = $shp.TextFrame.TextRange

# 2) Insert two new paragraphs after paragraph 4 ("Would require a lot of new code...")
#    and before paragraph 5 ("Generate probabilities by counting usage")
$para4 = $tr.Paragraphs(4, 1)
$para4.InsertAfter("`rWhen using textgenerationpipeline from transformers as installed dependency, it makes the typescript linter lag and doesn't work with webpack or web. It also costs a 700 kb of space.`rSetting xenova/transformers up turns out to be hard, because one of their dependencies requires FS.")

# Re-fetch text range after structural edit
$tr = $shp.TextFrame.TextRange

# paragraph 5 is now "When using textgenerationpipeline ... space."
$para5 = $tr.Paragraphs(5, 1)
# split into 3 runs: "When using " / "textgenerationpipeline" / " from transformers..."
$run1 = $para5.Characters(1, 11)
$run1.Text = "When using "
$run2 = $para5.Characters(12, 22)
$run2.Text = "textgenerationpipeline"
$run3len = $para5.Length - 33
$run3 = $para5.Characters(34, $run3len)
$run3.Text = " from transformers as installed dependency, it makes the typescript linter lag and doesn't work with webpack or web. It also costs a 700 kb of space."

# paragraph 6 is now "Setting xenova/transformers up turns out to be hard, because one of their dependencies requires FS."
$tr = $shp.TextFrame.TextRange
$para6 = $tr.Paragraphs(6, 1)
$run4 = $para6.Characters(1, 8)
$run4.Text = "Setting "
$run5 = $para6.Characters(9, 6)
$run5.Text = "xenova"
$run6len = $para6.Length - 14
$run6 = $para6.Characters(15, $run6len)
$run6.Text = "/transformers up turns out to be hard, because one of their dependencies requires FS."

# 3) Append a new paragraph after the last paragraph ("Would take less code and performance.")
$tr = $shp.TextFrame.TextRange
$count = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($count, 1)
$lastPara.InsertAfter("`rWouldn't be context aware.")

Write-Host "edit complete"
